$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 9).Value = 2566.6667   # I62
$ws.Cells.Item(62, 10).Value = 2750   # J62
$ws.Cells.Item(62, 11).Value = 2566.6667   # K62
$ws.Cells.Item(62, 12).Value = 2750   # L62
$ws.Cells.Item(62, 13).Value = -1942.6667   # M62
$ws.Cells.Item(62, 14).Value = -3998   # N62
$ws.Cells.Item(65, 9).Value = 2566.6667   # I65
$ws.Cells.Item(65, 10).Value = 2750   # J65
$ws.Cells.Item(65, 11).Value = 12833.3335   # K65
$ws.Cells.Item(65, 12).Value = 13750   # L65
$ws.Cells.Item(65, 13).Value = -9713.333500000001   # M65
$ws.Cells.Item(65, 14).Value = -19990   # N65
$ws.Cells.Item(74, 8).Value = 3531.6155   # H74
$ws.Cells.Item(74, 9).Value = 3416.375   # I74
$ws.Cells.Item(74, 10).Value = 3716   # J74
$ws.Cells.Item(74, 11).Value = 3416.375   # K74
$ws.Cells.Item(74, 12).Value = 3716   # L74
$ws.Cells.Item(74, 13).Value = -2480.375   # M74
$ws.Cells.Item(74, 14).Value = -5588   # N74
$ws.Cells.Item(77, 8).Value = 3531.6155   # H77
$ws.Cells.Item(77, 9).Value = 3416.375   # I77
$ws.Cells.Item(77, 10).Value = 3716   # J77
$ws.Cells.Item(77, 11).Value = 17081.875   # K77
$ws.Cells.Item(77, 12).Value = 18580   # L77
$ws.Cells.Item(77, 13).Value = -12401.875   # M77
$ws.Cells.Item(77, 14).Value = -27940   # N77
$ws.Cells.Item(132, 8).Value = 29412846   # H132
$ws.Cells.Item(132, 9).Value = 34483468   # I132
$ws.Cells.Item(132, 10).Value = 3249.2   # J132
$ws.Cells.Item(132, 11).Value = 103450404   # K132
$ws.Cells.Item(132, 12).Value = 9747.599999999999   # L132
$ws.Cells.Item(132, 13).Value = -103447874   # M132
$ws.Cells.Item(132, 14).Value = -14807.6   # N132
$ws.Cells.Item(137, 8).Value = 1025.6072   # H137
$ws.Cells.Item(137, 9).Value = 881.4545000000001   # I137
$ws.Cells.Item(137, 11).Value = 2644.3635   # K137
$ws.Cells.Item(137, 13).Value = -94.36350000000039   # M137

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 338929.28   # H32
$ws.Cells.Item(32, 9).Value = 2384.1519   # I32
$ws.Cells.Item(32, 11).Value = 2384.1519   # K32
$ws.Cells.Item(32, 13).Value = -2097.1519   # M32
$ws.Cells.Item(74, 8).Value = 670.46155   # H74
$ws.Cells.Item(74, 9).Value = 572.9091   # I74
$ws.Cells.Item(74, 10).Value = 1207   # J74
$ws.Cells.Item(74, 11).Value = 572.9091   # K74
$ws.Cells.Item(74, 12).Value = 1207   # L74
$ws.Cells.Item(74, 13).Value = 301.0909   # M74
$ws.Cells.Item(74, 14).Value = -2955   # N74
$ws.Cells.Item(77, 8).Value = 670.46155   # H77
$ws.Cells.Item(77, 9).Value = 572.9091   # I77
$ws.Cells.Item(77, 10).Value = 1207   # J77
$ws.Cells.Item(77, 11).Value = 2864.5455   # K77
$ws.Cells.Item(77, 12).Value = 6035   # L77
$ws.Cells.Item(77, 13).Value = 1503.4545   # M77
$ws.Cells.Item(77, 14).Value = -14771   # N77
$ws.Cells.Item(104, 8).Value = 16600   # H104
$ws.Cells.Item(104, 10).Value = 16600   # J104
$ws.Cells.Item(104, 12).Value = 16600   # L104
$ws.Cells.Item(104, 14).Value = -23588   # N104
$ws.Cells.Item(122, 8).Value = 12254.775   # H122
$ws.Cells.Item(122, 9).Value = 14200.714   # I122
$ws.Cells.Item(122, 10).Value = 1660.2222   # J122
$ws.Cells.Item(122, 11).Value = 42602.142   # K122
$ws.Cells.Item(122, 12).Value = 4980.6666   # L122
$ws.Cells.Item(122, 13).Value = -40152.142   # M122
$ws.Cells.Item(122, 14).Value = -9880.6666   # N122
$ws.Cells.Item(132, 8).Value = 19611116   # H132
$ws.Cells.Item(132, 9).Value = 25000892   # I132
$ws.Cells.Item(132, 10).Value = 11932.182   # J132
$ws.Cells.Item(132, 11).Value = 75002676   # K132
$ws.Cells.Item(132, 12).Value = 35796.546   # L132
$ws.Cells.Item(132, 13).Value = -75000146   # M132
$ws.Cells.Item(132, 14).Value = -40856.546   # N132

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 8967.954   # H134
$ws.Cells.Item(134, 9).Value = 2763.9473   # I134
$ws.Cells.Item(134, 10).Value = 48260   # J134
$ws.Cells.Item(134, 11).Value = 8291.841899999999   # K134
$ws.Cells.Item(134, 12).Value = 144780   # L134
$ws.Cells.Item(134, 13).Value = -5756.841899999999   # M134
$ws.Cells.Item(134, 14).Value = -149850   # N134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2971.45   # H31
$ws.Cells.Item(31, 9).Value = 2241.652   # I31
$ws.Cells.Item(31, 11).Value = 2241.652   # K31
$ws.Cells.Item(31, 13).Value = -1946.652   # M31
$ws.Cells.Item(34, 8).Value = 2971.45   # H34
$ws.Cells.Item(34, 9).Value = 2241.652   # I34
$ws.Cells.Item(34, 11).Value = 2241.652   # K34
$ws.Cells.Item(34, 13).Value = -2039.652   # M34
$ws.Cells.Item(107, 8).Value = 962.625   # H107
$ws.Cells.Item(107, 9).Value = 1170.3334   # I107
$ws.Cells.Item(107, 10).Value = 838   # J107
$ws.Cells.Item(107, 11).Value = 1170.3334   # K107
$ws.Cells.Item(107, 12).Value = 838   # L107
$ws.Cells.Item(107, 13).Value = 749.6666   # M107
$ws.Cells.Item(107, 14).Value = -4678   # N107
$ws.Cells.Item(134, 8).Value = 38566.258   # H134
$ws.Cells.Item(134, 9).Value = 51789.5   # I134
$ws.Cells.Item(134, 10).Value = 785.5714   # J134
$ws.Cells.Item(134, 11).Value = 155368.5   # K134
$ws.Cells.Item(134, 12).Value = 2356.7142   # L134
$ws.Cells.Item(134, 13).Value = -152833.5   # M134
$ws.Cells.Item(134, 14).Value = -7426.7142   # N134

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1035.7059   # H5
$ws.Cells.Item(5, 9).Value = 1150.2142   # I5
$ws.Cells.Item(5, 10).Value = 501.33334   # J5
$ws.Cells.Item(5, 11).Value = 3450.6426   # K5
$ws.Cells.Item(5, 12).Value = 1504.00002   # L5
$ws.Cells.Item(5, 13).Value = -3338.6426   # M5
$ws.Cells.Item(5, 14).Value = -1728.00002   # N5
$ws.Cells.Item(51, 8).Value = 1050.7273   # H51
$ws.Cells.Item(51, 9).Value = 150.42857   # I51
$ws.Cells.Item(51, 10).Value = 2626.25   # J51
$ws.Cells.Item(51, 11).Value = 451.28571   # K51
$ws.Cells.Item(51, 12).Value = 7878.75   # L51
$ws.Cells.Item(51, 13).Value = 8.714290000000005   # M51
$ws.Cells.Item(51, 14).Value = -8798.75   # N51
$ws.Cells.Item(129, 8).Value = 15875473   # H129
$ws.Cells.Item(129, 9).Value = 4900   # I129
$ws.Cells.Item(129, 10).Value = 18520568   # J129
$ws.Cells.Item(129, 11).Value = 14700   # K129
$ws.Cells.Item(129, 12).Value = 55561704   # L129
$ws.Cells.Item(129, 13).Value = -9700   # M129
$ws.Cells.Item(129, 14).Value = -55571704   # N129
$ws.Cells.Item(131, 8).Value = 6945399   # H131
$ws.Cells.Item(131, 9).Value = 686.6   # I131
$ws.Cells.Item(131, 10).Value = 8772955   # J131
$ws.Cells.Item(131, 11).Value = 2059.8   # K131
$ws.Cells.Item(131, 12).Value = 26318865   # L131
$ws.Cells.Item(131, 13).Value = 2980.2   # M131
$ws.Cells.Item(131, 14).Value = -26328945   # N131
$ws.Cells.Item(135, 8).Value = 1035.7059   # H135
$ws.Cells.Item(135, 9).Value = 1150.2142   # I135
$ws.Cells.Item(135, 10).Value = 501.33334   # J135
$ws.Cells.Item(135, 11).Value = 10351.9278   # K135
$ws.Cells.Item(135, 12).Value = 4512.00006   # L135
$ws.Cells.Item(135, 13).Value = -7816.927799999999   # M135
$ws.Cells.Item(135, 14).Value = -9582.00006   # N135

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 17052706   # H70
$ws.Cells.Item(70, 9).Value = 28131616   # I70
$ws.Cells.Item(70, 11).Value = 28131616   # K70
$ws.Cells.Item(70, 13).Value = -28131346   # M70
$ws.Cells.Item(73, 8).Value = 17052706   # H73
$ws.Cells.Item(73, 9).Value = 28131616   # I73
$ws.Cells.Item(73, 11).Value = 28131616   # K73
$ws.Cells.Item(73, 13).Value = -28130680   # M73
$ws.Cells.Item(132, 8).Value = 66326.23   # H132
$ws.Cells.Item(132, 9).Value = 46731.91   # I132
$ws.Cells.Item(132, 10).Value = 114223.445   # J132
$ws.Cells.Item(132, 11).Value = 140195.73   # K132
$ws.Cells.Item(132, 12).Value = 342670.335   # L132
$ws.Cells.Item(132, 13).Value = -137665.73   # M132
$ws.Cells.Item(132, 14).Value = -347730.335   # N132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(41, 8).Value = 1253824.5   # H41
$ws.Cells.Item(41, 10).Value = 5099   # J41
$ws.Cells.Item(41, 12).Value = 5099   # L41
$ws.Cells.Item(41, 14).Value = -5975   # N41
$ws.Cells.Item(55, 8).Value = 1079.5   # H55
$ws.Cells.Item(55, 9).Value = 1379   # I55
$ws.Cells.Item(55, 10).Value = 780   # J55
$ws.Cells.Item(55, 11).Value = 1379   # K55
$ws.Cells.Item(55, 12).Value = 780   # L55
$ws.Cells.Item(55, 13).Value = -1206   # M55
$ws.Cells.Item(55, 14).Value = -1126   # N55
$ws.Cells.Item(60, 8).Value = 0   # H60
$ws.Cells.Item(60, 10).Value = 0   # J60
$ws.Cells.Item(60, 12).Value = 0   # L60
$ws.Cells.Item(60, 14).ClearContents()   # N60 (remove cell)
$ws.Cells.Item(68, 8).Value = 2024.3611   # H68
$ws.Cells.Item(68, 9).Value = 1877.3914   # I68
$ws.Cells.Item(68, 11).Value = 1877.3914   # K68
$ws.Cells.Item(68, 13).Value = -1128.3914   # M68
$ws.Cells.Item(71, 8).Value = 2024.3611   # H71
$ws.Cells.Item(71, 9).Value = 1877.3914   # I71
$ws.Cells.Item(71, 11).Value = 9386.957   # K71
$ws.Cells.Item(71, 13).Value = -5642.957   # M71
$ws.Cells.Item(132, 8).Value = 4827.222   # H132
$ws.Cells.Item(132, 9).Value = 5016.4614   # I132
$ws.Cells.Item(132, 10).Value = 3597.1667   # J132
$ws.Cells.Item(132, 11).Value = 15049.3842   # K132
$ws.Cells.Item(132, 12).Value = 10791.5001   # L132
$ws.Cells.Item(132, 13).Value = -12519.3842   # M132
$ws.Cells.Item(132, 14).Value = -15851.5001   # N132
$ws.Cells.Item(136, 8).Value = 6591.6553   # H136
$ws.Cells.Item(136, 9).Value = 5746   # I136
$ws.Cells.Item(136, 11).Value = 17238   # K136
$ws.Cells.Item(136, 13).Value = -14688   # M136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(37, 8).Value = 3549.25   # H37
$ws.Cells.Item(37, 9).Value = 3400   # I37
$ws.Cells.Item(37, 11).Value = 3400   # K37
$ws.Cells.Item(37, 13).Value = -3197   # M37
$ws.Cells.Item(51, 8).Value = 1200   # H51
$ws.Cells.Item(51, 9).Value = 1200   # I51
$ws.Cells.Item(51, 11).Value = 1200   # K51
$ws.Cells.Item(51, 13).Value = -690   # M51
$ws.Cells.Item(54, 8).Value = 6270.125   # H54
$ws.Cells.Item(54, 10).Value = 6270.125   # J54
$ws.Cells.Item(54, 12).Value = 6270.125   # L54
$ws.Cells.Item(54, 14).Value = -7310.125   # N54
$ws.Cells.Item(123, 8).Value = 30000   # H123
$ws.Cells.Item(123, 10).Value = 30000   # J123
$ws.Cells.Item(123, 12).Value = 30000   # L123
$ws.Cells.Item(123, 14).Value = -39800   # N123
$ws.Cells.Item(126, 8).Value = 38463012   # H126
$ws.Cells.Item(126, 9).Value = 71429910   # I126
$ws.Cells.Item(126, 10).Value = 1632.0834   # J126
$ws.Cells.Item(126, 11).Value = 214289730   # K126
$ws.Cells.Item(126, 12).Value = 4896.2502   # L126
$ws.Cells.Item(126, 13).Value = -214287260   # M126
$ws.Cells.Item(126, 14).Value = -9836.2502   # N126
$ws.Cells.Item(132, 8).Value = 17331542   # H132
$ws.Cells.Item(132, 9).Value = 20835744   # I132
$ws.Cells.Item(132, 10).Value = 2313534   # J132
$ws.Cells.Item(132, 11).Value = 62507232   # K132
$ws.Cells.Item(132, 12).Value = 6940602   # L132
$ws.Cells.Item(132, 13).Value = -62504702   # M132
$ws.Cells.Item(132, 14).Value = -6945662   # N132
